$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tasks in rows 8 and 9 moved from "To Do" into "In Progress"
$ws.Range("D8").Value = "In Progress"
$ws.Range("D9").Value = "In Progress"

# Row 8's task picked up a recorded Time value now that work has started
$ws.Range("F8").Value = 0.17

# Re-enter the Overcost formula across the whole column so Excel stores it
# as a single shared formula (G2:G13), matching the recalculated values
$ws.Range("G2:G13").Formula = "=F2-E2"

# Column D is now noticeably wider text ("In Progress"), so it gets
# auto-fit like the sheet's other best-fit columns
$ws.Columns.Item(4).EntireColumn.AutoFit()

# Leave the selection where the user was last working
$ws.Range("C8:D9").Select()

$wb.Save()
